$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values to insert for rows 20-44 (A20:A44)
$values = @(
  '$ 27.354 CLP 23-10-20',
  '$ 27.354 CLP 23-10-20',
  '$ 27.381 CLP 28-10-20',
  '$ 27.381 CLP 28-10-20',
  '$ 27.381 CLP 28-10-20',
  '$ 27.386 CLP 29-10-20',
  '$ 27.386 CLP 29-10-20',
  '$ 27.386 CLP 29-10-20',
  '$ 27.386 CLP 29-10-20',
  '$ 27.386 CLP 29-10-20',
  '$ 27.386 CLP 29-10-20',
  '$ 27.386 CLP 29-10-20',
  '$ 27.386 CLP 29-10-20',
  '$ 27.386 CLP 29-10-20',
  '$ 27.386 CLP 29-10-20',
  '$ 27.386 CLP 29-10-20',
  '$ 27.386 CLP 29-10-20',
  '$ 27.386 CLP 29-10-20',
  '$ 27.386 CLP 29-10-20',
  '$ 27.386 CLP 29-10-20',
  '$ 27.386 CLP 29-10-20',
  '$ 27.386 CLP 29-10-20',
  '$ 27.386 CLP 29-10-20',
  '$ 27.386 CLP 29-10-20',
  '$ 27.386 CLP 29-10-20'
)

$row = 20
foreach ($val in $values) {
  $ws.Cells.Item($row, 1).Value = $val
  $row = $row + 1
}
